# Updated symbol list on Tue Dec 27 23:56:08 UTC 2022 with GitHub Actions
# Refreshes the "Price" (D), and for the re-ranked rows 18-24 also the
# "Coin" (B), "Link" (C) and "Volume(1h)" (E) columns, to match the latest
# coinranking.com snapshot. Numeric-looking prices are written with a
# leading "'" so Excel stores them as text (matching the source sheet's
# inlineStr cells) instead of silently coercing them to floats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.43"
$ws.Range("D3").Value = "'24.02"
$ws.Range("D4").Value = "'5.356"
$ws.Range("D5").Value = "'0.05811"
$ws.Range("D6").Value = "'3.376"
$ws.Range("D7").Value = "'6.474"
$ws.Range("D8").Value = "'0.8096"
$ws.Range("D9").Value = "'0.9240"
$ws.Range("D10").Value = "'0.1403"
$ws.Range("D11").Value = "'0.07405"
$ws.Range("D12").Value = "'0.03197"
$ws.Range("D13").Value = "'0.03028"
$ws.Range("D15").Value = "'3.860"
$ws.Range("D16").Value = "'0.001556"
$ws.Range("D17").Value = "'0.04696"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005987"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006176"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001254"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.004688"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00008796"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.597"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.150"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D28").Value = "'0.0002349"
$ws.Range("D40").Value = "'0.03847"
$ws.Range("D41").Value = "'0.006357"
$ws.Range("D42").Value = "'0.003498"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("D43").Value = "'0.1066"
$ws.Range("D44").Value = "'0.009067"
$ws.Range("D45").Value = "'0.00005267"
$ws.Range("D47").Value = "'0.6852"
$ws.Range("D48").Value = "'0.001843"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
